$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(6)

# The first paragraph's run currently reads "Sub/Factroy ". We need it to
# read "../sample/Sub/Factroy " as a single run (keeping its rPr/formatting
# intact). Because the new text shares "Sub/Factroy " as a trailing
# substring with the old text, a direct assignment would otherwise get
# smart-diffed into two runs (an inserted "../sample/" run plus the
# untouched original run). Stomping the paragraph text first with an
# unrelated placeholder collapses it back down to a single run (while
# preserving its rPr), so the following assignment produces one clean run.
$shape.TextFrame.TextRange.Paragraphs(1).Text = "X"
$shape.TextFrame.TextRange.Paragraphs(1).Text = "../sample/Sub/Factroy "
